# Apply updated absenteeism data (ETL consolidator refactor) to rows 2-11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  69375, "Vinicius Dias",            "Operações",              "Viagem de negócios", 8, 45078, 8341.07),
    @(3,  38900, "Augusto Moreira",          "P&D",                    "Doença",              3, 45089, 5377.27),
    @(4,  82807, "Stephany Rocha",           "Operações",              "Outros",              1, 45084, 12435.58),
    @(5,  10062, "Pedro Henrique Novaes",    "Jurídico",                "Problemas pessoais",  8, 45088, 2985.61),
    @(6,  86981, "Emanuel Santos",           "Atendimento ao Cliente", "Consulta médica",     8, 45081, 2733.46),
    @(7,  36339, "Rafaela Freitas",          "Engenharia",             "Problemas pessoais",  2, 45081, 4441.9),
    @(8,  62396, "Sra. Stella da Cunha",     "Financeiro",             "Problemas pessoais",  7, 45093, 5775.06),
    @(9,  69723, "Raul Rezende",             "TI",                     "Consulta médica",     8, 45088, 12396.95),
    @(10, 55356, "Maria Clara Pereira",      "TI",                     "Outros",              5, 45089, 9766.11),
    @(11, 52800, "Marcos Vinicius da Mata",  "TI",                     "Viagem de negócios",  1, 45078, 5855.8)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
